$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values for row 2 (B2:F2)
$ws.Range("B2").Value = 0.09494813238296555
$ws.Range("C2").Value = 0.5332219804412588
$ws.Range("D2").Value = 0.4006434568944247
$ws.Range("E2").Value = 0.6329640249606803
$ws.Range("F2").Value = 0.6398668578903548

# Shift previous rows 2-10 down to rows 3-11 (B:F columns)
$ws.Range("B3").Value = 0.1072470739083369
$ws.Range("C3").Value = 0.583924513487991
$ws.Range("D3").Value = 0.5944036395536997
$ws.Range("E3").Value = 0.7709757710549013
$ws.Range("F3").Value = 0.7806389925716273

$ws.Range("B4").Value = 0.573629067650584
$ws.Range("C4").Value = 0.8880024064378264
$ws.Range("D4").Value = 3.979434659499331
$ws.Range("E4").Value = 1.994852039500507
$ws.Range("F4").Value = 1.953538050233822

$ws.Range("B5").Value = 0.228737977167174
$ws.Range("C5").Value = 1.392066633737142
$ws.Range("D5").Value = 7.304711183963507
$ws.Range("E5").Value = 2.702722920308981
$ws.Range("F5").Value = 2.753551201305114

$ws.Range("B6").Value = 0.1376631994370348
$ws.Range("C6").Value = 1.217887797378426
$ws.Range("D6").Value = 7.521227355969471
$ws.Range("E6").Value = 2.742485616365101
$ws.Range("F6").Value = 2.80058721592678

$ws.Range("B7").Value = 0.256355043509169
$ws.Range("C7").Value = 1.359960095653506
$ws.Range("D7").Value = 7.005699073826143
$ws.Range("E7").Value = 2.646828115655821
$ws.Range("F7").Value = 2.693591441706439

$ws.Range("B8").Value = 0.09029628155329977
$ws.Range("C8").Value = 1.369716599985631
$ws.Range("D8").Value = 7.335564264182292
$ws.Range("E8").Value = 2.708424683128976
$ws.Range("F8").Value = 2.767756307023949

$ws.Range("B9").Value = 0.1622545047491004
$ws.Range("C9").Value = 1.487177404830128
$ws.Range("D9").Value = 8.202263652189526
$ws.Range("E9").Value = 2.863959436198342
$ws.Range("F9").Value = 2.92362286105626

$ws.Range("B10").Value = 0.08448834823307437
$ws.Range("C10").Value = 1.390260216460968
$ws.Range("D10").Value = 7.363914233767023
$ws.Range("E10").Value = 2.713653300214864
$ws.Range("F10").Value = 2.773296749308377

$ws.Range("B11").Value = 0.1037379453787874
$ws.Range("C11").Value = 1.518653016633174
$ws.Range("D11").Value = 7.970283297256951
$ws.Range("E11").Value = 2.823169016771216
$ws.Range("F11").Value = 2.88466950940459
